$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for each affected cell, keyed by row -> column letter
$newValues = @{
    2 = @{ "G"=30.29931333333333; "H"=90.89794000000001; "I"=0.2999359759116586; "J"=0.2999359759116585; "M"=0.01989833333333333; "N"=0.059695; "O"=0.5455534129646046; "P"=0.5455534129646046; "Q"=0.6029058364777778; "R"=5.4261525283; "S"=0.1636310953294748; "T"=0.1636310953294748 }
    3 = @{ "G"=30.29931333333333; "H"=90.89794000000001; "I"=0.2999359759116586; "J"=0.2999359759116585; "O"=0.1997148627777118; "P"=0.1997148627777118; "Q"=0.2207102980911111; "R"=1.98639268282; "S"=0.05990167227129596; "T"=0.05990167227129596 }
    4 = @{ "G"=30.29931333333333; "H"=90.89794000000001; "I"=0.2999359759116586; "J"=0.2999359759116585; "M"=0.009290999999999999; "N"=0.027873; "O"=0.2547317242576836; "P"=0.2547317242576836; "Q"=0.28151092018; "R"=2.53359828162; "S"=0.07640320831088784; "T"=0.07640320831088784 }
    5 = @{ "I"=0.1596558491596042; "J"=0.1596558491596042; "M"=0.01989833333333333; "N"=0.059695; "O"=0.5455534129646046; "P"=0.5455534129646046; "Q"=0.3209266344044444; "R"=2.88833970964; "S"=0.08710079340878418; "T"=0.08710079340878417 }
    6 = @{ "I"=0.1596558491596042; "J"=0.1596558491596042; "O"=0.1997148627777118; "P"=0.1997148627777118; "S"=0.03188564600656941; "T"=0.03188564600656941 }
    7 = @{ "I"=0.1596558491596042; "J"=0.1596558491596042; "M"=0.009290999999999999; "N"=0.027873; "O"=0.2547317242576836; "P"=0.2547317242576836; "Q"=0.149848196344; "R"=1.348633767096; "S"=0.04066940974425063; "T"=0.04066940974425063 }
    8 = @{ "G"=25.10878733333334; "H"=75.326362; "I"=0.2485544325684925; "J"=0.2485544325684925; "M"=0.01989833333333333; "N"=0.059695; "O"=0.5455534129646046; "P"=0.5455534129646046; "Q"=0.4996230199544445; "R"=4.49660717959; "S"=0.1355997189952217; "T"=0.1355997189952217 }
    9 = @{ "G"=25.10878733333334; "H"=75.326362; "I"=0.2485544325684925; "J"=0.2485544325684925; "O"=0.1997148627777118; "P"=0.1997148627777118; "Q"=0.1829007765317778; "R"=1.646106988786; "S"=0.0496400143932085; "T"=0.0496400143932085 }
    10 = @{ "G"=25.10878733333334; "H"=75.326362; "I"=0.2485544325684925; "J"=0.2485544325684925; "M"=0.009290999999999999; "N"=0.027873; "O"=0.2547317242576836; "P"=0.2547317242576836; "Q"=0.233285743114; "R"=2.099571688026; "S"=0.06331469918006224; "T"=0.06331469918006224 }
    11 = @{ "G"=29.482852; "H"=88.448556; "I"=0.2918537423602448; "J"=0.2918537423602447; "M"=0.01989833333333333; "N"=0.059695; "O"=0.5455534129646046; "P"=0.5455534129646046; "Q"=0.5866596167133333; "R"=5.27993655042; "S"=0.1592218052311239; "T"=0.1592218052311239 }
    12 = @{ "G"=29.482852; "H"=88.448556; "I"=0.2918537423602448; "J"=0.2918537423602447; "O"=0.1997148627777118; "P"=0.1997148627777118; "Q"=0.2147629215853333; "R"=1.932866294268; "S"=0.05828753010663792; "T"=0.05828753010663792 }
    13 = @{ "G"=29.482852; "H"=88.448556; "I"=0.2918537423602448; "J"=0.2918537423602447; "M"=0.009290999999999999; "N"=0.027873; "O"=0.2547317242576836; "P"=0.2547317242576836; "Q"=0.2739251779319999; "R"=2.465326601388; "S"=0.07434440702248289; "T"=0.07434440702248289 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowData[$colLetter]
    }
}
